# Insert a new data row right before the current row 296 ("Rabanito" /
# Vega Central Mapocho de Santiago weekly series), shifting all the
# subsequent rows down by one and populating the new row with a fresh
# weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 296:392 down to 297:393, leaving row 296 empty (but carrying
# the formatting of the row it was inserted above, matching column D's
# date style).
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new observation.
$ws.Range("A296").Value = 9
$ws.Range("B296").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C296").Value = "Metropolitana"
$ws.Range("D296").Value = 44985
$ws.Range("E296").Value = 13
$ws.Range("F296").Value = 300000001
$ws.Range("G296").Value = "Rabanito"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 7000
$ws.Range("K296").Value = 3000
$ws.Range("L296").Value = 3000
$ws.Range("M296").Value = 3000
$ws.Range("N296").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O296").Value = "Provincia de Chacabuco"
$ws.Range("P296").Value = 30
$ws.Range("Q296").Value = 100
$ws.Range("R296").Value = "Hortaliza"
